$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.028.41'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.103.23'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.62'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.32'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -7.42%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.093.50'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.64'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.159'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.11'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000218'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.596.06'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.974.69'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.43%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.099.23'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '506.31'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.68'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.57'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.708'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.26'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.96'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.38'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.75'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.96'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -8.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.35'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.53'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.47%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '59.54'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +12.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '530.90'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -12.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.22'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -7.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0412'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0794'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.060.71'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.73'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -9.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.10'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.255'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.06'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.91'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.15'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -7.00%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0512'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.48%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.107'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.66%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +52.94%  '
